$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------------
# Row 2: COO 154 -> 25004, Serie 789 -> 301
$ws.Range("A2").Value = 25004
$ws.Range("B2").Value = 301

# Row 3: COO 149 -> 25012 (Serie stays blank)
$ws.Range("A3").Value = 25012

# Rows 4-5: clear the old COO values (150 / 100); Serie column was already blank
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()

# --- Font / formatting -------------------------------------------------------
# Column A (rows 2-5) and B2 move to a "Courier New 8pt black" font and lose
# the centered alignment they used to share with the rest of the sheet.
# Build the format once on A2 (clearing the inherited underline first) and
# then fan it out with copy / paste-special so we don't re-derive the font
# table for every single cell.
$master = $ws.Range("A2")
$mf = $master.Font
$mf.Underline = -4142   # xlUnderlineStyleNone - clear any inherited underline
$mf.Name = "Courier New"
$mf.Family = 3          # Modern (monospace) family, matches Courier New
$mf.Size = 8
$mf.Color = 0
$master.HorizontalAlignment = 1     # xlGeneral
$master.VerticalAlignment = -4107   # xlBottom

$master.Copy()
foreach ($addr in @("A3", "A4", "A5", "B2")) {
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- Selection ----------------------------------------------------------------
$ws.Range("A3").Select() | Out-Null
